$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-11 Thursday" "2024-07-12 Friday"

Replace-Text "149×6=894" "568×2=1136"
Replace-Text "400×6=2400" "908×2=1816"
Replace-Text "106×7=742" "575×5=2875"
Replace-Text "507×4=2028" "792×7=5544"
Replace-Text "941×5=4705" "859×2=1718"

Replace-Text "206×4=824" "232×7=1624"
Replace-Text "215×7=1505" "661×2=1322"
Replace-Text "741×5=3705" "343×2=686"
Replace-Text "309×2=618" "452×6=2712"
Replace-Text "479×8=3832" "172×2=344"

Replace-Text "801×8=6408" "954×6=5724"
Replace-Text "563×3=1689" "297×3=891"
Replace-Text "415×9=3735" "660×9=5940"
Replace-Text "218×7=1526" "258×3=774"
Replace-Text "231×4=924" "817×9=7353"

Replace-Text "715×7=5005" "494×6=2964"
Replace-Text "327×2=654" "808×5=4040"
Replace-Text "999×6=5994" "667×4=2668"
Replace-Text "191×7=1337" "744×5=3720"
Replace-Text "695×3=2085" "137×9=1233"

Replace-Text "340×9=3060" "198×8=1584"
Replace-Text "663×6=3978" "141×2=282"
Replace-Text "535×7=3745" "193×2=386"
Replace-Text "161×4=644" "521×6=3126"
Replace-Text "989×6=5934" "260×5=1300"
